$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (the "LOW THREAT" / Nile Air NP-120 row) entirely.
# This shifts row 3 (Nile Air NP-110 / HIGH THREAT) up to become the new row 2,
# which also carries forward its own formatting already present in the sheet.
$ws.Rows.Item(2).Delete()
